$wb = $excel.ActiveWorkbook
$cards = $wb.Worksheets.Item("Cards")
$recipes = $wb.Worksheets.Item("Recipes")

# --- Cards sheet: text content updates ---------------------------------

# Row 2 - Bolt Rats
$cards.Range("G2").Value = "Strong Magnet`nfor 4 Metal"
$cards.Range("H2").Value = "6 HP. 2 ATK. `nInsta-slaughter: Use Tote Bag"

# Row 3 - Can of Lube
$cards.Range("G3").Value = "Fireball`nfor 1 Fabric + 1 Oil"

# Row 4 - Socket Wrench Puppy
$cards.Range("G4").Value = "Bandages`nfor 1 Fabric + 1 Duct Tape"

# Row 5 - Minor Ductwork
$cards.Range("G5").Value = "Tape Dispensor`nfor 2 Metal + 1 Oil"
$cards.Range("G5").WrapText = $true

# Row 6 - Rusty Toolbox
$cards.Range("G6").Value = "Strange Device`nfor 3 Metal + 1 Duct Tape"
$cards.Range("G6").WrapText = $true

# Row 7 - Tetanus Worms
$cards.Range("G7").Value = "Whacking Plank`nfor 3 Metal and 1 Duct tape"
$cards.Range("G7").WrapText = $true

# Row 8 - Curtain Ghost
$cards.Range("F8").Value = "1:2 for Metal:Fabric"
$cards.Range("G8").Value = "Tote Bag`nfor 3 Fabric"
$cards.Range("H8").Value = "18 HP. 2 ATK."

# Row 9 - Sprocket Bats
$cards.Range("H9").Value = "3 HP. 3 ATK.`nInsta-slaughter: Have Whacking Plank"

# Row 10 - Oil Nymph
$cards.Range("F10").Value = "1:1 Oil:Duct tape"
$cards.Range("G10").Value = "Fireball`nfor 1 Fabric + 1 Oil"
$cards.Range("H10").Value = "13 HP. 3 ATK. `nInsta-slaughter: Fireball"
$cards.Range("G10").WrapText = $true
$cards.Range("H10").WrapText = $true

# Row 11 - Petulant Tig Welder
$cards.Range("F11").Value = "4 Metal for +1 ATK"

# --- Cards sheet: row heights (auto-grown by the wrapped, multi-line text) ---
$cards.Rows.Item(3).RowHeight = 30
$cards.Rows.Item(6).RowHeight = 30
$cards.Rows.Item(8).RowHeight = 30
$cards.Rows.Item(9).RowHeight = 60
$cards.Rows.Item(10).RowHeight = 45

# --- Cards sheet: column widths for the re-flowed Build/Battle columns ---
$cards.Columns.Item(7).ColumnWidth = 24.5
$cards.Columns.Item(8).ColumnWidth = 34.6

# --- Selections -----------------------------------------------------------
$recipes.Activate()
$recipes.Range("F2").Select() | Out-Null

$cards.Activate()
$cards.Range("H3").Select() | Out-Null
